{"js": "// The resume's address line still shows a stale ZIP code (\"Fort Lauderdale, FL  32839\").\n// Strip the \" 32839\" (the extra space + zip code) so the line reads \"Fort Lauderdale, FL \".\nconst results = context.document.body.search(\" 32839\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items,text\");\nawait context.sync();\n\nfor (const result of results.items) {\n  result.insertText(\"\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The resume's address line still shows a stale ZIP code (\"Fort Lauderdale, FL  32839\").\n# Strip the \" 32839\" (the extra space + zip code) so the line reads \"Fort Lauderdale, FL \".\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \" 32839\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"\"\n$find.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 2) | Out-Null\n"}
